$d = $word.ActiveDocument
$d.Content.Find.Execute("Junior", $true, $false, $false, $false, $false, $true, 1, $false, "Trainee", 2)
